$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("delSubscriptionById")
$ws.Range("C4").Value = "testClient01"
$ws.Range("G4").Value = 200
$ws.Range("H4").Value = 107003
$ws.Range("I4").Value = "not exists"
$ws.Activate() | Out-Null
$ws.Range("I8").Select() | Out-Null
